# Auto-generated Excel COM-interop script applying the Diabolos_Profits market-price update.
# For each changed row, sets the recomputed price/profit columns (H-N) to their new values;
# where a cell is removed entirely in the target state, ClearContents() is used instead.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 191032.72
$ws.Range("J17").Value = 193937.19
$ws.Range("L17").Value = 581811.5700000001
$ws.Range("N17").Value = -582147.5700000001
$ws.Range("H19").Value = 3041.5334
$ws.Range("I19").Value = 2624.6667
$ws.Range("J19").Value = 3319.4443
$ws.Range("K19").Value = 2624.6667
$ws.Range("L19").Value = 3319.4443
$ws.Range("M19").Value = -2449.6667
$ws.Range("N19").Value = -3669.4443
$ws.Range("I28").Value = 39155.5
$ws.Range("J28").Value = 6144.375
$ws.Range("K28").Value = 39155.5
$ws.Range("L28").Value = 6144.375
$ws.Range("M28").Value = -38670.5
$ws.Range("N28").Value = -7114.375
$ws.Range("H62").Value = 832436.3
$ws.Range("I62").Value = 1146074.6
$ws.Range("J62").Value = 126750
$ws.Range("K62").Value = 1146074.6
$ws.Range("L62").Value = 126750
$ws.Range("M62").Value = -1145450.6
$ws.Range("N62").Value = -127998
$ws.Range("H65").Value = 832436.3
$ws.Range("I65").Value = 1146074.6
$ws.Range("J65").Value = 126750
$ws.Range("K65").Value = 5730373
$ws.Range("L65").Value = 633750
$ws.Range("M65").Value = -5727253
$ws.Range("N65").Value = -639990
$ws.Range("H86").Value = 28724284
$ws.Range("I86").Value = 4997.5
$ws.Range("K86").Value = 4997.5
$ws.Range("M86").Value = -3874.5
$ws.Range("H88").Value = 1002.7
$ws.Range("I88").Value = 669.6667
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 669.6667
$ws.Range("L88").Value = 4000
$ws.Range("M88").Value = -263.6667
$ws.Range("N88").Value = -4812
$ws.Range("H89").Value = 28724284
$ws.Range("I89").Value = 4997.5
$ws.Range("K89").Value = 24987.5
$ws.Range("M89").Value = -19371.5
$ws.Range("H91").Value = 1002.7
$ws.Range("I91").Value = 669.6667
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 669.6667
$ws.Range("L91").Value = 4000
$ws.Range("M91").Value = 734.3333
$ws.Range("N91").Value = -6808
$ws.Range("H98").Value = 586.871
$ws.Range("I98").Value = 506.65518
$ws.Range("K98").Value = 506.65518
$ws.Range("M98").Value = 991.34482
$ws.Range("H116").Value = 73045300
$ws.Range("J116").Value = 66671268
$ws.Range("L116").Value = 66671268
$ws.Range("N116").Value = -66678152
$ws.Range("H122").Value = 586.871
$ws.Range("I122").Value = 506.65518
$ws.Range("K122").Value = 1519.96554
$ws.Range("M122").Value = 930.0344600000001
$ws.Range("H125").Value = 1698.4
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1698.4
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 15285.6
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -20205.6
$ws.Range("H135").Value = 1670.6471
$ws.Range("I135").Value = 738.75
$ws.Range("J135").Value = 3907.2
$ws.Range("K135").Value = 6648.75
$ws.Range("L135").Value = 35164.8
$ws.Range("M135").Value = -4113.75
$ws.Range("N135").Value = -40234.8
$ws.Range("H137").Value = 28573224
$ws.Range("I137").Value = 47620136
$ws.Range("K137").Value = 142860408
$ws.Range("M137").Value = -142857858
$ws.Range("H138").Value = 2259.7874
$ws.Range("I138").Value = 978.0571
$ws.Range("K138").Value = 2934.1713
$ws.Range("M138").Value = 2205.8287
$ws.Range("H141").Value = 843.2273
$ws.Range("I141").Value = 870.0476
$ws.Range("K141").Value = 2610.1428
$ws.Range("M141").Value = 2569.8572

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H32").Value = 3885.5789
$ws.Range("I32").Value = 2254.7144
$ws.Range("K32").Value = 2254.7144
$ws.Range("M32").Value = -1967.7144
$ws.Range("H45").Value = 2469.3333
$ws.Range("I45").Value = 2046.2858
$ws.Range("K45").Value = 2046.2858
$ws.Range("M45").Value = -1669.2858
$ws.Range("H132").Value = 1716.4
$ws.Range("I132").Value = 1573.7778
$ws.Range("K132").Value = 4721.3334
$ws.Range("M132").Value = -2191.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2746.5715
$ws.Range("I99").Value = 2506.3333
$ws.Range("J99").Value = 2926.75
$ws.Range("K99").Value = 2506.3333
$ws.Range("L99").Value = 2926.75
$ws.Range("M99").Value = -1008.3333
$ws.Range("N99").Value = -5922.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 752.41174
$ws.Range("I22").Value = 690.9167
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 690.9167
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -340.9167
$ws.Range("N22").Value = -1600
$ws.Range("H31").Value = 3344.6667
$ws.Range("I31").Value = 1749.7142
$ws.Range("K31").Value = 1749.7142
$ws.Range("M31").Value = -1454.7142
$ws.Range("H34").Value = 3344.6667
$ws.Range("I34").Value = 1749.7142
$ws.Range("K34").Value = 1749.7142
$ws.Range("M34").Value = -1547.7142
$ws.Range("H82").Value = 25000
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25722
$ws.Range("H85").Value = 25000
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27496
$ws.Range("H134").Value = 577.5789
$ws.Range("I134").Value = 527.05884
$ws.Range("K134").Value = 1581.17652
$ws.Range("M134").Value = 953.82348
$ws.Range("H140").Value = 167244
$ws.Range("J140").Value = 191666.5
$ws.Range("L140").Value = 191666.5
$ws.Range("N140").Value = -202026.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 14401015
$ws.Range("I4").Value = 26548444
$ws.Range("K4").Value = 79645332
$ws.Range("M4").Value = -79645220
$ws.Range("H9").Value = 108894.7
$ws.Range("J9").Value = 170500.5
$ws.Range("L9").Value = 511501.5
$ws.Range("N9").Value = -511949.5
$ws.Range("H23").Value = 73.75
$ws.Range("I23").Value = 41.666668
$ws.Range("J23").Value = 93
$ws.Range("K23").Value = 125.000004
$ws.Range("L23").Value = 279
$ws.Range("M23").Value = 109.999996
$ws.Range("N23").Value = -749
$ws.Range("H29").Value = 198.625
$ws.Range("I29").Value = 331.75
$ws.Range("K29").Value = 995.25
$ws.Range("M29").Value = -718.25
$ws.Range("H133").Value = 4378.2
$ws.Range("I133").Value = 4378.2
$ws.Range("K133").Value = 13134.6
$ws.Range("M133").Value = -8074.599999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 32.666668
$ws.Range("I2").Value = 36.76923
$ws.Range("J2").Value = 6
$ws.Range("K2").Value = 36.76923
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 76.23077000000001
$ws.Range("N2").Value = -232
$ws.Range("H25").Value = 2499
$ws.Range("I25").Value = 2498
$ws.Range("J25").Value = 2499.5
$ws.Range("K25").Value = 2498
$ws.Range("L25").Value = 2499.5
$ws.Range("M25").Value = -1969
$ws.Range("N25").Value = -3557.5
$ws.Range("H43").Value = 2728.3333
$ws.Range("I43").Value = 1361
$ws.Range("J43").Value = 4437.5
$ws.Range("K43").Value = 1361
$ws.Range("L43").Value = 4437.5
$ws.Range("M43").Value = -1210
$ws.Range("N43").Value = -4739.5
$ws.Range("H113").Value = 1882.9048
$ws.Range("I113").Value = 1653.1875
$ws.Range("K113").Value = 1653.1875
$ws.Range("M113").Value = 516.8125
$ws.Range("H121").Value = 20000
$ws.Range("J121").Value = 20000
$ws.Range("L121").Value = 20000
$ws.Range("N121").Value = -23494

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 643.25
$ws.Range("J22").Value = 689.4
$ws.Range("L22").Value = 689.4
$ws.Range("N22").Value = -1279.4
$ws.Range("H27").Value = 643.25
$ws.Range("J27").Value = 689.4
$ws.Range("L27").Value = 689.4
$ws.Range("N27").Value = -903.4
$ws.Range("H46").Value = 3129.1428
$ws.Range("I46").Value = 1644.5454
$ws.Range("K46").Value = 1644.5454
$ws.Range("M46").Value = -1456.5454
$ws.Range("H55").Value = 428.48
$ws.Range("I55").Value = 279
$ws.Range("J55").Value = 901.8333
$ws.Range("K55").Value = 279
$ws.Range("L55").Value = 901.8333
$ws.Range("M55").Value = -106
$ws.Range("N55").Value = -1247.8333
$ws.Range("H132").Value = 6854.1665
$ws.Range("I132").Value = 3190.4614
$ws.Range("J132").Value = 16379.8
$ws.Range("K132").Value = 9571.3842
$ws.Range("L132").Value = 49139.39999999999
$ws.Range("M132").Value = -7041.3842
$ws.Range("N132").Value = -54199.39999999999
$ws.Range("H136").Value = 2230.8845
$ws.Range("I136").Value = 2080.16
$ws.Range("K136").Value = 6240.48
$ws.Range("M136").Value = -3690.48

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 12507990
$ws.Range("I81").Value = 3147.6667
$ws.Range("J81").Value = 20010896
$ws.Range("K81").Value = 6295.3334
$ws.Range("L81").Value = 40021792
$ws.Range("M81").Value = -5234.3334
$ws.Range("N81").Value = -40023914
$ws.Range("H84").Value = 12507990
$ws.Range("I84").Value = 3147.6667
$ws.Range("J84").Value = 20010896
$ws.Range("K84").Value = 31476.667
$ws.Range("L84").Value = 200108960
$ws.Range("M84").Value = -26172.667
$ws.Range("N84").Value = -200119568
$ws.Range("H119").Value = 43474.5
$ws.Range("J119").Value = 41669.4
$ws.Range("L119").Value = 41669.4
$ws.Range("N119").Value = -51345.4
$ws.Range("H122").Value = 1546.3334
$ws.Range("I122").Value = 1434.9286
$ws.Range("K122").Value = 4304.7858
$ws.Range("M122").Value = -1854.7858

